$wb = $excel.ActiveWorkbook

$oldGuid = "ffde527a-175c-46a2-9f59-448d398e30c5"
$newGuid = "67d2aa2d-b14b-4712-8fcc-2614d98023f2"

# Sheet "Overview"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-11-29 03:14:44"

# Sheet "zh-cn"
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.d90d391f81e2a2a0af158aea9367903a4c41db2a.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-11-29 03:14:31"

# Sheet "de-de"
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.d90d391f81e2a2a0af158aea9367903a4c41db2a.de-de.xlf"
